$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the two orphaned/duplicate hidden chart-tracking defined names.
#    (_xlchart.v1.10 / .11 duplicate .8 / .9 - Excel prunes them on save.)
# ---------------------------------------------------------------------------
$wb.Names("_xlchart.v1.10").Delete()
$wb.Names("_xlchart.v1.11").Delete()

# ---------------------------------------------------------------------------
# 2. Add the "Mean increase" / "Median increase" headers + formulas to each
#    of the four summary blocks (Low, Medium, High, All).
# ---------------------------------------------------------------------------
function Add-IncreaseBlock {
    param(
        [int]$HeaderRow,
        [int]$FormulaRow,
        [string]$MeanFormula,
        [string]$MedianFormula
    )

    $ws.Cells.Item($HeaderRow, 4).Value = "Mean increase"
    $ws.Cells.Item($HeaderRow, 4).Font.Bold = $true

    $ws.Cells.Item($HeaderRow, 6).Value = "Median increase"
    $ws.Cells.Item($HeaderRow, 6).Font.Bold = $true

    $ws.Cells.Item($FormulaRow, 4).Formula = $MeanFormula
    $ws.Cells.Item($FormulaRow, 4).ClearFormats()

    $ws.Cells.Item($FormulaRow, 6).Formula = $MedianFormula
    $ws.Cells.Item($FormulaRow, 6).ClearFormats()
}

# Block 1 - Low (rows 2:31 / B2:B31, mean E3:E4, median E10)
Add-IncreaseBlock -HeaderRow 18 -FormulaRow 19 `
    -MeanFormula "= ((E3 / 114.202998) * 100) - 100" `
    -MedianFormula "= ((E10 / 113.658804) * 100) - 100"

# Block 2 - Medium (rows 34:63 / B34:B63, mean E35:E36, median E42)
Add-IncreaseBlock -HeaderRow 50 -FormulaRow 51 `
    -MeanFormula "= ((E35 / 114.202998) * 100) - 100" `
    -MedianFormula "= ((E42 / 113.658804) * 100) - 100"

# Block 3 - High (rows 66:95 / B66:B95, mean E67:E68, median E74)
Add-IncreaseBlock -HeaderRow 82 -FormulaRow 83 `
    -MeanFormula "= ((E67 / 114.202998) * 100) - 100" `
    -MedianFormula "= ((E74 / 113.658804) * 100) - 100"

# Block 4 - All (rows 98:187, averages the three groups above)
Add-IncreaseBlock -HeaderRow 113 -FormulaRow 114 `
    -MeanFormula "= (D19 + D51 + D83) / 3" `
    -MedianFormula "= (F19 + F51 + F83) / 3"

# ---------------------------------------------------------------------------
# 3. Update the view: scroll position + active selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 87
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F104").Select()
